$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "access time" (column C) and "size" (column D) values for rows 2-25
# to reflect new measurements taken after moving testing to Mac (Safari).
# All "size" values become 0.

$data = @{
    2  = 0.1803960800170898
    3  = 0.1174867153167725
    4  = 0.3515560626983643
    5  = 0.1403241157531738
    6  = 0.04364013671875
    7  = 0.0446479320526123
    8  = 0.1361439228057861
    9  = 0.1383850574493408
    10 = 0.03995084762573242
    11 = 0.04306912422180176
    12 = 0.1463279724121094
    13 = 0.1357040405273438
    14 = 0.05635976791381836
    15 = 0.04310894012451172
    16 = 0.151137113571167
    17 = 0.1682040691375732
    18 = 0.04387784004211426
    19 = 0.04345989227294922
    20 = 0.1562318801879883
    21 = 0.1395237445831299
    22 = 0.06763029098510742
    23 = 0.04022407531738281
    24 = 0.1517479419708252
    25 = 0.1604301929473877
}

foreach ($row in $data.Keys) {
    $ws.Range("C$row").Value = $data[$row]
    $ws.Range("D$row").Value = 0
}
